# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.267.64"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "'2.064.08"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'233.76"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'56.55"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'2.366.95"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "'14.59"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "'0.775"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "'2.063.02"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "'37.255.64"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  +5.52%  "
$ws.Range("D20").Value = "'69.39"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "'226.31"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").Value = "'166.26"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").Value = "'19.02"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "'4.45"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("E39").Value = "  -4.76%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").Value = "'1.465.98"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'96.04"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("D46").Value = "'4.28"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'14.95"
$ws.Range("E48").Value = "  -7.24%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.11"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "'2.252.95"
$ws.Range("E51").Value = "  -0.20%  "
